$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.007.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.72"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5102"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06359"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.67"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07770"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.266"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.86"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5440"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7723"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.24"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.975.84"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.65"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.12%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.898"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.047"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.890"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.40"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1203"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.831"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.60"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.234"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04865"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.279"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.182"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.531"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.375"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9090"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.586"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.127.07"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5468"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01562"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.51%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.523"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8095"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.65%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.30"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.49%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.464"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.91%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₈123"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.773.79"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4533"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.22%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.88"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05295"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.33%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.11%  "
